$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply header-style formatting (matching row 20 J-column style) to the newly
# populated J21:J23 header cells before writing their values.
$ws.Cells.Item(20, 10).Copy() | Out-Null
$ws.Cells.Item(21, 10).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22, 10).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(23, 10).PasteSpecial(-4122) | Out-Null

# Row 1
$ws.Cells.Item(1, 1).Value = "negative"
$ws.Cells.Item(1, 10).Value = "positive"

# Row 2
$ws.Cells.Item(2, 1).Value = "name"
$ws.Cells.Item(2, 2).Value = "anchor score"
$ws.Cells.Item(2, 3).Value = "type occurences"
$ws.Cells.Item(2, 4).Value = "total occurences"
$ws.Cells.Item(2, 5).Value = "+%"
$ws.Cells.Item(2, 6).Value = "-%"
$ws.Cells.Item(2, 7).Value = "both"
$ws.Cells.Item(2, 8).Value = "normal"
$ws.Cells.Item(2, 10).Value = "name"
$ws.Cells.Item(2, 11).Value = "anchor score"
$ws.Cells.Item(2, 12).Value = "type occurences"
$ws.Cells.Item(2, 13).Value = "total occurences"
$ws.Cells.Item(2, 14).Value = "+%"
$ws.Cells.Item(2, 15).Value = "-%"
$ws.Cells.Item(2, 16).Value = "both"
$ws.Cells.Item(2, 17).Value = "normal"

# Row 3
$ws.Cells.Item(3, 1).Value = "poorly"
$ws.Cells.Item(3, 2).Value = 0.9565217391304348
$ws.Cells.Item(3, 3).Value = 44
$ws.Cells.Item(3, 4).Value = 44
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = $false
$ws.Cells.Item(3, 8).Value = 2
$ws.Cells.Item(3, 10).Value = "awesome"
$ws.Cells.Item(3, 11).Value = 0.9076923076923077
$ws.Cells.Item(3, 12).Value = 59
$ws.Cells.Item(3, 13).Value = 59
$ws.Cells.Item(3, 14).Value = 1
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = $false
$ws.Cells.Item(3, 17).Value = 6

# Row 4
$ws.Cells.Item(4, 1).Value = "disappointing"
$ws.Cells.Item(4, 2).Value = 0.7954545454545454
$ws.Cells.Item(4, 3).Value = 35
$ws.Cells.Item(4, 4).Value = 35
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = $false
$ws.Cells.Item(4, 8).Value = 9
$ws.Cells.Item(4, 10).Value = "wonderful"
$ws.Cells.Item(4, 11).Value = 0.8928571428571429
$ws.Cells.Item(4, 12).Value = 50
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 1
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = $false
$ws.Cells.Item(4, 17).Value = 6

# Row 5
$ws.Cells.Item(5, 1).Value = "poor"
$ws.Cells.Item(5, 2).Value = 0.7464788732394366
$ws.Cells.Item(5, 3).Value = 53
$ws.Cells.Item(5, 4).Value = 53
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = $false
$ws.Cells.Item(5, 8).Value = 18
$ws.Cells.Item(5, 10).Value = "favorite"
$ws.Cells.Item(5, 11).Value = 0.8064516129032258
$ws.Cells.Item(5, 12).Value = 75
$ws.Cells.Item(5, 13).Value = 75
$ws.Cells.Item(5, 14).Value = 1
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = $false
$ws.Cells.Item(5, 17).Value = 18

# Row 6
$ws.Cells.Item(6, 1).Value = "disappointed"
$ws.Cells.Item(6, 2).Value = 0.6935483870967742
$ws.Cells.Item(6, 3).Value = 129
$ws.Cells.Item(6, 4).Value = 129
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = $false
$ws.Cells.Item(6, 8).Value = 57
$ws.Cells.Item(6, 10).Value = "excellent"
$ws.Cells.Item(6, 11).Value = 0.765625
$ws.Cells.Item(6, 12).Value = 49
$ws.Cells.Item(6, 13).Value = 49
$ws.Cells.Item(6, 14).Value = 1
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = $false
$ws.Cells.Item(6, 17).Value = 15

# Row 7
$ws.Cells.Item(7, 1).Value = "however"
$ws.Cells.Item(7, 2).Value = 0.671875
$ws.Cells.Item(7, 3).Value = 43
$ws.Cells.Item(7, 4).Value = 43
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = $false
$ws.Cells.Item(7, 8).Value = 21
$ws.Cells.Item(7, 10).Value = "classic"
$ws.Cells.Item(7, 11).Value = 0.660377358490566
$ws.Cells.Item(7, 12).Value = 35
$ws.Cells.Item(7, 13).Value = 35
$ws.Cells.Item(7, 14).Value = 1
$ws.Cells.Item(7, 15).Value = 0
$ws.Cells.Item(7, 16).Value = $false
$ws.Cells.Item(7, 17).Value = 18

# Row 8
$ws.Cells.Item(8, 1).Value = "junk"
$ws.Cells.Item(8, 2).Value = 0.6363636363636364
$ws.Cells.Item(8, 3).Value = 35
$ws.Cells.Item(8, 4).Value = 35
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = 20
$ws.Cells.Item(8, 10).Value = "love"
$ws.Cells.Item(8, 11).Value = 0.5451936872309899
$ws.Cells.Item(8, 12).Value = 380
$ws.Cells.Item(8, 13).Value = 380
$ws.Cells.Item(8, 14).Value = 1
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = $false
$ws.Cells.Item(8, 17).Value = 317

# Row 9
$ws.Cells.Item(9, 1).Value = "waste"
$ws.Cells.Item(9, 2).Value = 0.6351351351351351
$ws.Cells.Item(9, 3).Value = 94
$ws.Cells.Item(9, 4).Value = 94
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = $false
$ws.Cells.Item(9, 8).Value = 54
$ws.Cells.Item(9, 10).Value = "thank"
$ws.Cells.Item(9, 11).Value = 0.5362318840579711
$ws.Cells.Item(9, 12).Value = 37
$ws.Cells.Item(9, 13).Value = 37
$ws.Cells.Item(9, 14).Value = 1
$ws.Cells.Item(9, 15).Value = 0
$ws.Cells.Item(9, 16).Value = $false
$ws.Cells.Item(9, 17).Value = 32

# Row 10
$ws.Cells.Item(10, 1).Value = "broke"
$ws.Cells.Item(10, 2).Value = 0.6067961165048543
$ws.Cells.Item(10, 3).Value = 125
$ws.Cells.Item(10, 4).Value = 125
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = $false
$ws.Cells.Item(10, 8).Value = 81
$ws.Cells.Item(10, 10).Value = "loves"
$ws.Cells.Item(10, 11).Value = 0.4937759336099585
$ws.Cells.Item(10, 12).Value = 238
$ws.Cells.Item(10, 13).Value = 238
$ws.Cells.Item(10, 14).Value = 1
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = $false
$ws.Cells.Item(10, 17).Value = 244

# Row 11
$ws.Cells.Item(11, 1).Value = "smaller"
$ws.Cells.Item(11, 2).Value = 0.5462184873949579
$ws.Cells.Item(11, 3).Value = 65
$ws.Cells.Item(11, 4).Value = 65
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = $false
$ws.Cells.Item(11, 8).Value = 54
$ws.Cells.Item(11, 10).Value = "great"
$ws.Cells.Item(11, 11).Value = 0.4467213114754098
$ws.Cells.Item(11, 12).Value = 545
$ws.Cells.Item(11, 13).Value = 545
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = $false
$ws.Cells.Item(11, 17).Value = 675

# Row 12
$ws.Cells.Item(12, 1).Value = "guess"
$ws.Cells.Item(12, 2).Value = 0.5370370370370371
$ws.Cells.Item(12, 3).Value = 29
$ws.Cells.Item(12, 4).Value = 29
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = $false
$ws.Cells.Item(12, 8).Value = 25
$ws.Cells.Item(12, 10).Value = "friends"
$ws.Cells.Item(12, 11).Value = 0.3544973544973545
$ws.Cells.Item(12, 12).Value = 67
$ws.Cells.Item(12, 13).Value = 67
$ws.Cells.Item(12, 14).Value = 1
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = $false
$ws.Cells.Item(12, 17).Value = 122

# Row 13
$ws.Cells.Item(13, 1).Value = "small"
$ws.Cells.Item(13, 2).Value = 0.4695652173913044
$ws.Cells.Item(13, 3).Value = 162
$ws.Cells.Item(13, 4).Value = 162
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = $false
$ws.Cells.Item(13, 8).Value = 183
$ws.Cells.Item(13, 10).Value = "loved"
$ws.Cells.Item(13, 11).Value = 0.3486238532110092
$ws.Cells.Item(13, 12).Value = 114
$ws.Cells.Item(13, 13).Value = 114
$ws.Cells.Item(13, 14).Value = 1
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = $false
$ws.Cells.Item(13, 17).Value = 213

# Row 14
$ws.Cells.Item(14, 1).Value = "broken"
$ws.Cells.Item(14, 2).Value = 0.4096385542168675
$ws.Cells.Item(14, 3).Value = 34
$ws.Cells.Item(14, 4).Value = 34
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = $false
$ws.Cells.Item(14, 8).Value = 49
$ws.Cells.Item(14, 10).Value = "perfect"
$ws.Cells.Item(14, 11).Value = 0.3192771084337349
$ws.Cells.Item(14, 12).Value = 53
$ws.Cells.Item(14, 13).Value = 53
$ws.Cells.Item(14, 14).Value = 1
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = $false
$ws.Cells.Item(14, 17).Value = 113

# Row 15
$ws.Cells.Item(15, 1).Value = "apart"
$ws.Cells.Item(15, 2).Value = 0.4
$ws.Cells.Item(15, 3).Value = 38
$ws.Cells.Item(15, 4).Value = 38
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = $false
$ws.Cells.Item(15, 8).Value = 57
$ws.Cells.Item(15, 10).Value = "best"
$ws.Cells.Item(15, 11).Value = 0.3166666666666667
$ws.Cells.Item(15, 12).Value = 38
$ws.Cells.Item(15, 13).Value = 38
$ws.Cells.Item(15, 14).Value = 1
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = $false
$ws.Cells.Item(15, 17).Value = 82

# Row 16
$ws.Cells.Item(16, 1).Value = "cheap"
$ws.Cells.Item(16, 2).Value = 0.3838862559241706
$ws.Cells.Item(16, 3).Value = 81
$ws.Cells.Item(16, 4).Value = 81
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = $false
$ws.Cells.Item(16, 8).Value = 130
$ws.Cells.Item(16, 10).Value = "happy"
$ws.Cells.Item(16, 11).Value = 0.2307692307692308
$ws.Cells.Item(16, 12).Value = 33
$ws.Cells.Item(16, 13).Value = 33
$ws.Cells.Item(16, 14).Value = 1
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = $false
$ws.Cells.Item(16, 17).Value = 110

# Row 17
$ws.Cells.Item(17, 1).Value = "plastic"
$ws.Cells.Item(17, 2).Value = 0.3779527559055118
$ws.Cells.Item(17, 3).Value = 48
$ws.Cells.Item(17, 4).Value = 48
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = $false
$ws.Cells.Item(17, 8).Value = 79
$ws.Cells.Item(17, 10).Value = "enjoy"
$ws.Cells.Item(17, 11).Value = 0.1935483870967742
$ws.Cells.Item(17, 12).Value = 36
$ws.Cells.Item(17, 13).Value = 36
$ws.Cells.Item(17, 14).Value = 1
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = $false
$ws.Cells.Item(17, 17).Value = 150

# Row 18
$ws.Cells.Item(18, 1).Value = "ok"
$ws.Cells.Item(18, 2).Value = 0.359375
$ws.Cells.Item(18, 3).Value = 46
$ws.Cells.Item(18, 4).Value = 46
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = $false
$ws.Cells.Item(18, 8).Value = 82
$ws.Cells.Item(18, 10).Value = "christmas"
$ws.Cells.Item(18, 11).Value = 0.1526104417670683
$ws.Cells.Item(18, 12).Value = 38
$ws.Cells.Item(18, 13).Value = 38
$ws.Cells.Item(18, 14).Value = 1
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = $false
$ws.Cells.Item(18, 17).Value = 211

# Row 19
$ws.Cells.Item(19, 1).Value = "thought"
$ws.Cells.Item(19, 2).Value = 0.2524752475247525
$ws.Cells.Item(19, 3).Value = 51
$ws.Cells.Item(19, 4).Value = 51
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = $false
$ws.Cells.Item(19, 8).Value = 151
$ws.Cells.Item(19, 10).Value = "fun"
$ws.Cells.Item(19, 11).Value = 0.1498685363716039
$ws.Cells.Item(19, 12).Value = 171
$ws.Cells.Item(19, 13).Value = 171
$ws.Cells.Item(19, 14).Value = 1
$ws.Cells.Item(19, 15).Value = 0
$ws.Cells.Item(19, 16).Value = $false
$ws.Cells.Item(19, 17).Value = 970

# Row 20
$ws.Cells.Item(20, 1).Value = "size"
$ws.Cells.Item(20, 2).Value = 0.2216494845360825
$ws.Cells.Item(20, 3).Value = 43
$ws.Cells.Item(20, 4).Value = 43
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = $false
$ws.Cells.Item(20, 8).Value = 151
$ws.Cells.Item(20, 10).Value = "easy"
$ws.Cells.Item(20, 11).Value = 0.0855614973262032
$ws.Cells.Item(20, 12).Value = 32
$ws.Cells.Item(20, 13).Value = 32
$ws.Cells.Item(20, 14).Value = 1
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = $false
$ws.Cells.Item(20, 17).Value = 342

# Row 21
$ws.Cells.Item(21, 1).Value = "hard"
$ws.Cells.Item(21, 2).Value = 0.175
$ws.Cells.Item(21, 3).Value = 35
$ws.Cells.Item(21, 4).Value = 35
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = $false
$ws.Cells.Item(21, 8).Value = 165
$ws.Cells.Item(21, 10).Value = "family"
$ws.Cells.Item(21, 11).Value = 0.0807799442896936
$ws.Cells.Item(21, 12).Value = 29
$ws.Cells.Item(21, 13).Value = 29
$ws.Cells.Item(21, 14).Value = 1
$ws.Cells.Item(21, 15).Value = 0
$ws.Cells.Item(21, 16).Value = $false
$ws.Cells.Item(21, 17).Value = 330

# Row 22
$ws.Cells.Item(22, 1).Value = "money"
$ws.Cells.Item(22, 2).Value = 0.1645569620253164
$ws.Cells.Item(22, 3).Value = 52
$ws.Cells.Item(22, 4).Value = 52
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = $false
$ws.Cells.Item(22, 8).Value = 264
$ws.Cells.Item(22, 10).Value = "game"
$ws.Cells.Item(22, 11).Value = 0.07667316439246263
$ws.Cells.Item(22, 12).Value = 118
$ws.Cells.Item(22, 13).Value = 120
$ws.Cells.Item(22, 14).Value = 0.98
$ws.Cells.Item(22, 15).Value = 0.02000000000000002
$ws.Cells.Item(22, 16).Value = $true
$ws.Cells.Item(22, 17).Value = 1421

# Row 23
$ws.Cells.Item(23, 1).Value = "item"
$ws.Cells.Item(23, 2).Value = 0.1630434782608696
$ws.Cells.Item(23, 3).Value = 45
$ws.Cells.Item(23, 4).Value = 45
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = $false
$ws.Cells.Item(23, 8).Value = 231
$ws.Cells.Item(23, 10).Value = "play"
$ws.Cells.Item(23, 11).Value = 0.05059920106524634
$ws.Cells.Item(23, 12).Value = 38
$ws.Cells.Item(23, 13).Value = 39
$ws.Cells.Item(23, 14).Value = 0.97
$ws.Cells.Item(23, 15).Value = 0.03000000000000003
$ws.Cells.Item(23, 16).Value = $true
$ws.Cells.Item(23, 17).Value = 713

# Row 24
$ws.Cells.Item(24, 1).Value = "price"
$ws.Cells.Item(24, 2).Value = 0.1556195965417868
$ws.Cells.Item(24, 3).Value = 54
$ws.Cells.Item(24, 4).Value = 55
$ws.Cells.Item(24, 5).Value = 0.02
$ws.Cells.Item(24, 6).Value = 0.98
$ws.Cells.Item(24, 7).Value = $true
$ws.Cells.Item(24, 8).Value = 293

# Row 25
$ws.Cells.Item(25, 1).Value = "would"
$ws.Cells.Item(25, 2).Value = 0.150297619047619
$ws.Cells.Item(25, 3).Value = 101
$ws.Cells.Item(25, 4).Value = 103
$ws.Cells.Item(25, 5).Value = 0.02
$ws.Cells.Item(25, 6).Value = 0.98
$ws.Cells.Item(25, 7).Value = $true
$ws.Cells.Item(25, 8).Value = 571

# Row 26
$ws.Cells.Item(26, 1).Value = "work"
$ws.Cells.Item(26, 2).Value = 0.1455696202531646
$ws.Cells.Item(26, 3).Value = 46
$ws.Cells.Item(26, 4).Value = 46
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = $false
$ws.Cells.Item(26, 8).Value = 270

# Row 27
$ws.Cells.Item(27, 1).Value = "product"
$ws.Cells.Item(27, 2).Value = 0.1258278145695364
$ws.Cells.Item(27, 3).Value = 57
$ws.Cells.Item(27, 4).Value = 58
$ws.Cells.Item(27, 5).Value = 0.02
$ws.Cells.Item(27, 6).Value = 0.98
$ws.Cells.Item(27, 7).Value = $true
$ws.Cells.Item(27, 8).Value = 396

# Row 28
$ws.Cells.Item(28, 1).Value = "use"
$ws.Cells.Item(28, 2).Value = 0.0958904109589041
$ws.Cells.Item(28, 3).Value = 35
$ws.Cells.Item(28, 4).Value = 35
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = $false
$ws.Cells.Item(28, 8).Value = 330

# Row 29
$ws.Cells.Item(29, 1).Value = "little"
$ws.Cells.Item(29, 2).Value = 0.08482142857142858
$ws.Cells.Item(29, 3).Value = 38
$ws.Cells.Item(29, 4).Value = 39
$ws.Cells.Item(29, 5).Value = 0.03
$ws.Cells.Item(29, 6).Value = 0.97
$ws.Cells.Item(29, 7).Value = $true
$ws.Cells.Item(29, 8).Value = 410

# Row 30
$ws.Cells.Item(30, 1).Value = "like"
$ws.Cells.Item(30, 2).Value = 0.066006600660066
$ws.Cells.Item(30, 3).Value = 40
$ws.Cells.Item(30, 4).Value = 42
$ws.Cells.Item(30, 5).Value = 0.05
$ws.Cells.Item(30, 6).Value = 0.95
$ws.Cells.Item(30, 7).Value = $true
$ws.Cells.Item(30, 8).Value = 566

# Row 31
$ws.Cells.Item(31, 1).Value = "one"
$ws.Cells.Item(31, 2).Value = 0.04309252217997465
$ws.Cells.Item(31, 3).Value = 34
$ws.Cells.Item(31, 4).Value = 39
$ws.Cells.Item(31, 5).Value = 0.13
$ws.Cells.Item(31, 6).Value = 0.87
$ws.Cells.Item(31, 7).Value = $true
$ws.Cells.Item(31, 8).Value = 755

# Row 32 ("one") no longer appears in the updated (min-count=5) results; remove it
# so the table ends at row 31 and the dimension shrinks to A1:Q31.
$ws.Rows.Item(32).Delete()
